# Script automates data entry into the bot's spreadsheet. Values are only
# poked in after the relevant UI elements/cells are ready/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name in A2 ("Rittmang" -> "Akhil")
$ws.Range("A2").Value = "Akhil"

# Update the phone/id number in B2
$ws.Range("B2").Value = 9921164006

# Move the selection/active cell to B3, ready for the next piece of data
$ws.Range("B3").Select()

# Maximize the Excel window (best effort - mirrors the workbookView window
# size/position captured when the workbook was last saved)
$excel.ActiveWindow.WindowState = -4137
$excel.ActiveWindow.Left = -108
$excel.ActiveWindow.Top = -108
$excel.ActiveWindow.Width = 23256
$excel.ActiveWindow.Height = 13176
